$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.517.91'
$ws.Range("E2").Value = '  +0.36%  '

# Row 3
$ws.Range("D3").Value = '1.907.75'
$ws.Range("E3").Value = '  +0.00%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.56%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.90'
$ws.Range("E5").Value = '  -0.46%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.40%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4844'
$ws.Range("E7").Value = '  +3.68%  '

# Row 8
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08153'
$ws.Range("E9").Value = '  +1.81%  '

# Row 10
$ws.Range("E10").Value = '  +0.44%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.44'
$ws.Range("E11").Value = '  +5.26%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.029'
$ws.Range("E12").Value = '  +1.55%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.847.44'
$ws.Range("E13").Value = '  -3.21%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.098'
$ws.Range("E14").Value = '  -0.36%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.43'
$ws.Range("E15").Value = '  +1.45%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.007'
$ws.Range("E16").Value = '  +0.51%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06752'
$ws.Range("E17").Value = '  +2.38%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001045'
$ws.Range("E18").Value = '  +1.56%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.71'
$ws.Range("E19").Value = '  +0.13%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.39%  '

# Row 21
$ws.Range("D21").Value = '29.530.38'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.597'
$ws.Range("E22").Value = '  +1.35%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.82'
$ws.Range("E23").Value = '  +2.72%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.166'
$ws.Range("E24").Value = '  -2.08%  '

# Row 25
$ws.Range("D25").Value = '2.138.20'
$ws.Range("E25").Value = '  +0.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.22'
$ws.Range("E26").Value = '  +0.46%  '

# Row 27
$ws.Range("E27").Value = '  +1.80%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.302'
$ws.Range("E28").Value = '  +10.55%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.105'
$ws.Range("E29").Value = '  -1.01%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.05'
$ws.Range("E30").Value = '  +1.89%  '

# Row 31
$ws.Range("E31").Value = '  -3.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09557'
$ws.Range("E32").Value = '  +1.06%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.537'
$ws.Range("E33").Value = '  +2.95%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.395'
$ws.Range("E34").Value = '  -1.48%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.555'
$ws.Range("E35").Value = '  -0.60%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02267'
$ws.Range("E36").Value = '  +0.69%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06117'
$ws.Range("E37").Value = '  +0.65%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.175'
$ws.Range("E38").Value = '  +0.22%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5969'
$ws.Range("E39").Value = '  +1.87%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.939'
$ws.Range("E40").Value = '  -5.15%  '

# Row 41
$ws.Range("E41").Value = '  +3.66%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1856'
$ws.Range("E42").Value = '  +1.24%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.437'
$ws.Range("E43").Value = '  +2.35%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.281'
$ws.Range("E44").Value = '  -1.58%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07721'
$ws.Range("E45").Value = '  +0.35%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.47'
$ws.Range("E46").Value = '  +2.78%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5577'
$ws.Range("E47").Value = '  +0.69%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.960'
$ws.Range("E48").Value = '  +1.96%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.99'
$ws.Range("E49").Value = '  +1.54%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.57'
$ws.Range("E50").Value = '  +1.81%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.055'
$ws.Range("E51").Value = '  +2.62%  '
